# Update the fixed "datetimeFigureOut" date text shown on the slide
# master and every slide layout's Date Placeholder (2023-10-07 -> 2023-10-11).
$p = $ppt.ActivePresentation

$newDate = "2023-10-11"

$sm = $p.SlideMaster

# Slide master's own Date Placeholder shape.
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# Each slide layout's own Date Placeholder shape.
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Resize/reposition the "Group 31" group shape on slide 2 (the group's
# bounding box was dragged to a new size/position; its child coordinate
# space - chOff/chExt - stays the same).
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Name -eq "Group 31") {
        $sh.Width = 137.45003
        $sh.Height = 35.63882
        $sh.Top = 825.61112
    }
}
